$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values (recalculated figures) ---
$ws.Range("G3").Value = 1.00002300652815
$ws.Range("G4").Value = 1.00002300652815
$ws.Range("G9").Value = 2.08680537211842
$ws.Range("G10").Value = 2.08680537211842
$ws.Range("G11").Value = 2.89530734800615
$ws.Range("G12").Value = 2.89530734800615
$ws.Range("G13").Value = 5.35608965517241
$ws.Range("L13").Value = 2.4216
$ws.Range("M13").Value = 11.89165
$ws.Range("G14").Value = 5.35608965517241
$ws.Range("L14").Value = 2.4216
$ws.Range("M14").Value = 11.89165
$ws.Range("G20").Value = 1.3716092134247
$ws.Range("G21").Value = 1.3716092134247
$ws.Range("G26").Value = 2.65671912394739
$ws.Range("G27").Value = 2.65671912394739
$ws.Range("G28").Value = 2.57550090869474
$ws.Range("G29").Value = 2.57550090869474
$ws.Range("F30").Value = 2.9246
$ws.Range("G30").Value = 6.5252275862069
$ws.Range("I30").Value = 25.59138
$ws.Range("M30").Value = 14.02292
$ws.Range("F31").Value = 2.9246
$ws.Range("G31").Value = 6.5252275862069
$ws.Range("I31").Value = 25.59138
$ws.Range("M31").Value = 14.02292
$ws.Range("G37").Value = 1.68438890631055
$ws.Range("G38").Value = 1.68438890631055
$ws.Range("G39").Value = 4910.26357783953
$ws.Range("I39").Value = 23695
$ws.Range("N39").Value = 9861.370269999999
$ws.Range("G40").Value = 4910.26357783953
$ws.Range("I40").Value = 23695
$ws.Range("N40").Value = 9861.370269999999
$ws.Range("G41").Value = 4910.26357783953
$ws.Range("I41").Value = 23695
$ws.Range("N41").Value = 9861.370269999999
$ws.Range("G42").Value = 4910.26357783953
$ws.Range("I42").Value = 23695
$ws.Range("N42").Value = 9861.370269999999
$ws.Range("G43").Value = 2.74432947834345
$ws.Range("G44").Value = 2.74432947834345
$ws.Range("F47").Value = 4.40435
$ws.Range("G47").Value = 7.211615
$ws.Range("I47").Value = 25.21565
$ws.Range("M47").Value = 13.78013
$ws.Range("F48").Value = 4.40435
$ws.Range("G48").Value = 7.211615
$ws.Range("I48").Value = 25.21565
$ws.Range("M48").Value = 13.78013
$ws.Range("G56").Value = 6678.69691117286
$ws.Range("N56").Value = 20417
$ws.Range("G57").Value = 6678.69691117286
$ws.Range("N57").Value = 20417
$ws.Range("G58").Value = 6678.69691117286
$ws.Range("N58").Value = 20417
$ws.Range("G59").Value = 6678.69691117286
$ws.Range("N59").Value = 20417
$ws.Range("G60").Value = 2.40651998582205
$ws.Range("G61").Value = 2.40651998582205
$ws.Range("F64").Value = 4.40435
$ws.Range("G64").Value = 7.047965
$ws.Range("M64").Value = 13.5211
$ws.Range("F65").Value = 4.40435
$ws.Range("G65").Value = 7.047965
$ws.Range("M65").Value = 13.5211
$ws.Range("G73").Value = 3902.54544122432
$ws.Range("I73").Value = 20500
$ws.Range("N73").Value = 13973
$ws.Range("G74").Value = 3902.54544122432
$ws.Range("I74").Value = 20500
$ws.Range("N74").Value = 13973
$ws.Range("G75").Value = 3902.54544122432
$ws.Range("I75").Value = 20500
$ws.Range("N75").Value = 13973
$ws.Range("G76").Value = 3902.54544122432
$ws.Range("I76").Value = 20500
$ws.Range("N76").Value = 13973
$ws.Range("G77").Value = 1.7722380311883
$ws.Range("G78").Value = 1.7722380311883
$ws.Range("G81").Value = 6.01724833333333
$ws.Range("L81").Value = 2.4058
$ws.Range("M81").Value = 10.56351
$ws.Range("N81").Value = 14.67554
$ws.Range("G82").Value = 6.01724833333333
$ws.Range("L82").Value = 2.4058
$ws.Range("M82").Value = 10.56351
$ws.Range("N82").Value = 14.67554
$ws.Range("G89").Value = 3907.37877455766
$ws.Range("I89").Value = 20500
$ws.Range("N89").Value = 13973
$ws.Range("G90").Value = 3907.37877455766
$ws.Range("I90").Value = 20500
$ws.Range("N90").Value = 13973
$ws.Range("G91").Value = 3907.37877455766
$ws.Range("I91").Value = 20500
$ws.Range("N91").Value = 13973
$ws.Range("G92").Value = 3907.37877455766
$ws.Range("I92").Value = 20500
$ws.Range("N92").Value = 13973
$ws.Range("G93").Value = 1.66970967354946
$ws.Range("G94").Value = 1.66970967354946
$ws.Range("G97").Value = 5.77134166666667
$ws.Range("N97").Value = 13.5467
$ws.Range("G98").Value = 5.77134166666667
$ws.Range("N98").Value = 13.5467
$ws.Range("G105").Value = 3778.12877455766
$ws.Range("I105").Value = 20500
$ws.Range("N105").Value = 12435.0444
$ws.Range("G106").Value = 3778.12877455766
$ws.Range("I106").Value = 20500
$ws.Range("N106").Value = 12435.0444
$ws.Range("G107").Value = 3778.12877455766
$ws.Range("I107").Value = 20500
$ws.Range("N107").Value = 12435.0444
$ws.Range("G108").Value = 3778.12877455766
$ws.Range("I108").Value = 20500
$ws.Range("N108").Value = 12435.0444
$ws.Range("G109").Value = 1.24855529830633
$ws.Range("G110").Value = 1.24855529830633
$ws.Range("G113").Value = 5.03520833333333
$ws.Range("I113").Value = 14.44595
$ws.Range("M113").Value = 7.26468
$ws.Range("N113").Value = 10.68783
$ws.Range("G114").Value = 5.03520833333333
$ws.Range("I114").Value = 14.44595
$ws.Range("M114").Value = 7.26468
$ws.Range("N114").Value = 10.68783
$ws.Range("G121").Value = 3602.49003056081
$ws.Range("I121").Value = 20750
$ws.Range("G122").Value = 3602.49003056081
$ws.Range("I122").Value = 20750
$ws.Range("G123").Value = 3602.49003056081
$ws.Range("I123").Value = 20750
$ws.Range("G124").Value = 3602.49003056081
$ws.Range("I124").Value = 20750
$ws.Range("G125").Value = 1.21549519444993
$ws.Range("G126").Value = 1.21549519444993
$ws.Range("G127").Value = 2.6883512651708
$ws.Range("G128").Value = 2.6883512651708
$ws.Range("G137").Value = 1792.43918310318
$ws.Range("H137").Value = 23171.9118030876
$ws.Range("G138").Value = 1792.43918310318
$ws.Range("H138").Value = 23171.9118030876
$ws.Range("G139").Value = 1792.43918310318
$ws.Range("H139").Value = 23171.9118030876
$ws.Range("G140").Value = 1792.43918310318
$ws.Range("H140").Value = 23171.9118030876
$ws.Range("G141").Value = 0.965263099102447
$ws.Range("G142").Value = 0.965263099102447
$ws.Range("G143").Value = 2.39962245161148
$ws.Range("G144").Value = 2.39962245161148
$ws.Range("G153").Value = 1326.02848344268
$ws.Range("H153").Value = 16855.2749688859
$ws.Range("G154").Value = 1326.02848344268
$ws.Range("H154").Value = 16855.2749688859
$ws.Range("G155").Value = 1326.02848344268
$ws.Range("H155").Value = 16855.2749688859
$ws.Range("G156").Value = 1326.02848344268
$ws.Range("H156").Value = 16855.2749688859
$ws.Range("G157").Value = 1.1541599192152
$ws.Range("G158").Value = 1.1541599192152
$ws.Range("G159").Value = 2.17423663181168
$ws.Range("G160").Value = 2.17423663181168

# --- Append new rows 167-182 ---
$arr = New-Object 'object[,]' 1,21
$arr[0,0] = 'Tutaenui Stream at d/s Marton STP'
$arr[0,1] = 'DRP (95th Percentile)'
$arr[0,2] = 'D'
$arr[0,3] = '2019 - 2023'
$arr[0,4] = 'Impact'
$arr[0,5] = 0.9615
$arr[0,6] = 1.64091379310345
$arr[0,7] = 5.37
$arr[0,8] = 4.498
$arr[0,11] = 2.34
$arr[0,12] = 3.9516
$arr[0,13] = 4.1886
$arr[0,14] = 1803578.705
$arr[0,15] = 5557699.998
$arr[0,16] = 'Rangitikei District'
$arr[0,17] = 'Rangitīkei-Turakina'
$arr[0,18] = 'Coastal Rangitikei'
$arr[0,19] = 'Rang_4d'
$arr[0,20] = 'mg/L'
$ws.Range("A167:U167").Value = $arr

$arr = New-Object 'object[,]' 1,21
$arr[0,0] = 'Tutaenui Stream at d/s Marton STP'
$arr[0,1] = 'DRP (Median)'
$arr[0,2] = 'D'
$arr[0,3] = '2019 - 2023'
$arr[0,4] = 'Impact'
$arr[0,5] = 0.9615
$arr[0,6] = 1.64091379310345
$arr[0,7] = 5.37
$arr[0,8] = 4.498
$arr[0,11] = 2.34
$arr[0,12] = 3.9516
$arr[0,13] = 4.1886
$arr[0,14] = 1803578.705
$arr[0,15] = 5557699.998
$arr[0,16] = 'Rangitikei District'
$arr[0,17] = 'Rangitīkei-Turakina'
$arr[0,18] = 'Coastal Rangitikei'
$arr[0,19] = 'Rang_4d'
$arr[0,20] = 'mg/L'
$ws.Range("A168:U168").Value = $arr

$arr = New-Object 'object[,]' 1,21
$arr[0,0] = 'Tutaenui Stream at d/s Marton STP'
$arr[0,1] = 'E coli (>260)'
$arr[0,2] = 'E'
$arr[0,3] = '2019 - 2023'
$arr[0,4] = 'Impact'
$arr[0,5] = 490
$arr[0,6] = 1252.90779378751
$arr[0,7] = 16855.2749688859
$arr[0,8] = 4080
$arr[0,9] = 44.8275862068966
$arr[0,10] = 72.4137931034483
$arr[0,11] = 899
$arr[0,12] = 1726.24
$arr[0,13] = 3601.38
$arr[0,14] = 1803578.705
$arr[0,15] = 5557699.998
$arr[0,16] = 'Rangitikei District'
$arr[0,17] = 'Rangitīkei-Turakina'
$arr[0,18] = 'Coastal Rangitikei'
$arr[0,19] = 'Rang_4d'
$arr[0,20] = '% exceedances over 260/100 mL'
$ws.Range("A169:U169").Value = $arr

$arr = New-Object 'object[,]' 1,21
$arr[0,0] = 'Tutaenui Stream at d/s Marton STP'
$arr[0,1] = 'E coli (>540)'
$arr[0,2] = 'E'
$arr[0,3] = '2019 - 2023'
$arr[0,4] = 'Impact'
$arr[0,5] = 490
$arr[0,6] = 1252.90779378751
$arr[0,7] = 16855.2749688859
$arr[0,8] = 4080
$arr[0,9] = 44.8275862068966
$arr[0,10] = 72.4137931034483
$arr[0,11] = 899
$arr[0,12] = 1726.24
$arr[0,13] = 3601.38
$arr[0,14] = 1803578.705
$arr[0,15] = 5557699.998
$arr[0,16] = 'Rangitikei District'
$arr[0,17] = 'Rangitīkei-Turakina'
$arr[0,18] = 'Coastal Rangitikei'
$arr[0,19] = 'Rang_4d'
$arr[0,20] = '% exceedances over 540/100 mL'
$ws.Range("A170:U170").Value = $arr

$arr = New-Object 'object[,]' 1,21
$arr[0,0] = 'Tutaenui Stream at d/s Marton STP'
$arr[0,1] = 'E coli (Median)'
$arr[0,2] = 'E'
$arr[0,3] = '2019 - 2023'
$arr[0,4] = 'Impact'
$arr[0,5] = 490
$arr[0,6] = 1252.90779378751
$arr[0,7] = 16855.2749688859
$arr[0,8] = 4080
$arr[0,9] = 44.8275862068966
$arr[0,10] = 72.4137931034483
$arr[0,11] = 899
$arr[0,12] = 1726.24
$arr[0,13] = 3601.38
$arr[0,14] = 1803578.705
$arr[0,15] = 5557699.998
$arr[0,16] = 'Rangitikei District'
$arr[0,17] = 'Rangitīkei-Turakina'
$arr[0,18] = 'Coastal Rangitikei'
$arr[0,19] = 'Rang_4d'
$arr[0,20] = 'E. coli/100 mL'
$ws.Range("A171:U171").Value = $arr

$arr = New-Object 'object[,]' 1,21
$arr[0,0] = 'Tutaenui Stream at d/s Marton STP'
$arr[0,1] = 'E coli (95th Percentile)'
$arr[0,2] = 'E'
$arr[0,3] = '2019 - 2023'
$arr[0,4] = 'Impact'
$arr[0,5] = 490
$arr[0,6] = 1252.90779378751
$arr[0,7] = 16855.2749688859
$arr[0,8] = 4080
$arr[0,9] = 44.8275862068966
$arr[0,10] = 72.4137931034483
$arr[0,11] = 899
$arr[0,12] = 1726.24
$arr[0,13] = 3601.38
$arr[0,14] = 1803578.705
$arr[0,15] = 5557699.998
$arr[0,16] = 'Rangitikei District'
$arr[0,17] = 'Rangitīkei-Turakina'
$arr[0,18] = 'Coastal Rangitikei'
$arr[0,19] = 'Rang_4d'
$arr[0,20] = 'E. coli/100 mL'
$ws.Range("A172:U172").Value = $arr

$arr = New-Object 'object[,]' 1,21
$arr[0,0] = 'Tutaenui Stream at d/s Marton STP'
$arr[0,1] = 'Ammoniacal-N (95th Percentile)'
$arr[0,2] = 'D'
$arr[0,3] = '2019 - 2023'
$arr[0,4] = 'Impact'
$arr[0,5] = 0.34953
$arr[0,6] = 0.846009912745724
$arr[0,7] = 5.62710303552482
$arr[0,8] = 3.67071
$arr[0,11] = 0.21639
$arr[0,12] = 1.76941
$arr[0,13] = 3.31961
$arr[0,14] = 1803578.705
$arr[0,15] = 5557699.998
$arr[0,16] = 'Rangitikei District'
$arr[0,17] = 'Rangitīkei-Turakina'
$arr[0,18] = 'Coastal Rangitikei'
$arr[0,19] = 'Rang_4d'
$arr[0,20] = 'mg NH4-N/L'
$ws.Range("A173:U173").Value = $arr

$arr = New-Object 'object[,]' 1,21
$arr[0,0] = 'Tutaenui Stream at d/s Marton STP'
$arr[0,1] = 'Ammoniacal-N (Median)'
$arr[0,2] = 'C'
$arr[0,3] = '2019 - 2023'
$arr[0,4] = 'Impact'
$arr[0,5] = 0.34953
$arr[0,6] = 0.846009912745724
$arr[0,7] = 5.62710303552482
$arr[0,8] = 3.67071
$arr[0,11] = 0.21639
$arr[0,12] = 1.76941
$arr[0,13] = 3.31961
$arr[0,14] = 1803578.705
$arr[0,15] = 5557699.998
$arr[0,16] = 'Rangitikei District'
$arr[0,17] = 'Rangitīkei-Turakina'
$arr[0,18] = 'Coastal Rangitikei'
$arr[0,19] = 'Rang_4d'
$arr[0,20] = 'mg NH4-N/L'
$ws.Range("A174:U174").Value = $arr

$arr = New-Object 'object[,]' 1,21
$arr[0,0] = 'Tutaenui Stream at d/s Marton STP'
$arr[0,1] = 'Nitrate-N (95th Percentile)'
$arr[0,2] = 'C'
$arr[0,3] = '2019 - 2023'
$arr[0,4] = 'Impact'
$arr[0,5] = 1.34
$arr[0,6] = 1.80206421801858
$arr[0,7] = 7.11
$arr[0,8] = 5.652
$arr[0,11] = 0.497
$arr[0,12] = 3.1296
$arr[0,13] = 5.1952
$arr[0,14] = 1803578.705
$arr[0,15] = 5557699.998
$arr[0,16] = 'Rangitikei District'
$arr[0,17] = 'Rangitīkei-Turakina'
$arr[0,18] = 'Coastal Rangitikei'
$arr[0,19] = 'Rang_4d'
$arr[0,20] = 'mg NO3-N/L'
$ws.Range("A175:U175").Value = $arr

$arr = New-Object 'object[,]' 1,21
$arr[0,0] = 'Tutaenui Stream at d/s Marton STP'
$arr[0,1] = 'Nitrate-N (Median)'
$arr[0,2] = 'B'
$arr[0,3] = '2019 - 2023'
$arr[0,4] = 'Impact'
$arr[0,5] = 1.34
$arr[0,6] = 1.80206421801858
$arr[0,7] = 7.11
$arr[0,8] = 5.652
$arr[0,11] = 0.497
$arr[0,12] = 3.1296
$arr[0,13] = 5.1952
$arr[0,14] = 1803578.705
$arr[0,15] = 5557699.998
$arr[0,16] = 'Rangitikei District'
$arr[0,17] = 'Rangitīkei-Turakina'
$arr[0,18] = 'Coastal Rangitikei'
$arr[0,19] = 'Rang_4d'
$arr[0,20] = 'mg NO3-N/L'
$ws.Range("A176:U176").Value = $arr

$arr = New-Object 'object[,]' 1,21
$arr[0,0] = 'Tutaenui Stream at d/s Marton STP'
$arr[0,1] = 'Soluble Inorganic Nitrogen (95th Percentile)'
$arr[0,3] = '2019 - 2023'
$arr[0,4] = 'Impact'
$arr[0,5] = 2.475
$arr[0,6] = 3.4061724137931
$arr[0,7] = 13.45
$arr[0,8] = 9.92
$arr[0,11] = 1.688
$arr[0,12] = 5.85572
$arr[0,13] = 7.7998
$arr[0,14] = 1803578.705
$arr[0,15] = 5557699.998
$arr[0,16] = 'Rangitikei District'
$arr[0,17] = 'Rangitīkei-Turakina'
$arr[0,18] = 'Coastal Rangitikei'
$arr[0,19] = 'Rang_4d'
$arr[0,20] = 'g/m3'
$ws.Range("A177:U177").Value = $arr

$arr = New-Object 'object[,]' 1,21
$arr[0,0] = 'Tutaenui Stream at d/s Marton STP'
$arr[0,1] = 'Soluble Inorganic Nitrogen (Median)'
$arr[0,3] = '2019 - 2023'
$arr[0,4] = 'Impact'
$arr[0,5] = 2.475
$arr[0,6] = 3.4061724137931
$arr[0,7] = 13.45
$arr[0,8] = 9.92
$arr[0,11] = 1.688
$arr[0,12] = 5.85572
$arr[0,13] = 7.7998
$arr[0,14] = 1803578.705
$arr[0,15] = 5557699.998
$arr[0,16] = 'Rangitikei District'
$arr[0,17] = 'Rangitīkei-Turakina'
$arr[0,18] = 'Coastal Rangitikei'
$arr[0,19] = 'Rang_4d'
$arr[0,20] = 'g/m3'
$ws.Range("A178:U178").Value = $arr

$arr = New-Object 'object[,]' 1,21
$arr[0,0] = 'Tutaenui Stream at d/s Marton STP'
$arr[0,1] = 'Total Nitrogen (95th Percentile)'
$arr[0,3] = '2019 - 2023'
$arr[0,4] = 'Impact'
$arr[0,5] = 5.53
$arr[0,6] = 6.61672413793103
$arr[0,7] = 17
$arr[0,8] = 14.9
$arr[0,11] = 6.75
$arr[0,12] = 10.564
$arr[0,13] = 12.4
$arr[0,14] = 1803578.705
$arr[0,15] = 5557699.998
$arr[0,16] = 'Rangitikei District'
$arr[0,17] = 'Rangitīkei-Turakina'
$arr[0,18] = 'Coastal Rangitikei'
$arr[0,19] = 'Rang_4d'
$arr[0,20] = 'g/m3'
$ws.Range("A179:U179").Value = $arr

$arr = New-Object 'object[,]' 1,21
$arr[0,0] = 'Tutaenui Stream at d/s Marton STP'
$arr[0,1] = 'Total Nitrogen (Median)'
$arr[0,3] = '2019 - 2023'
$arr[0,4] = 'Impact'
$arr[0,5] = 5.53
$arr[0,6] = 6.61672413793103
$arr[0,7] = 17
$arr[0,8] = 14.9
$arr[0,11] = 6.75
$arr[0,12] = 10.564
$arr[0,13] = 12.4
$arr[0,14] = 1803578.705
$arr[0,15] = 5557699.998
$arr[0,16] = 'Rangitikei District'
$arr[0,17] = 'Rangitīkei-Turakina'
$arr[0,18] = 'Coastal Rangitikei'
$arr[0,19] = 'Rang_4d'
$arr[0,20] = 'g/m3'
$ws.Range("A180:U180").Value = $arr

$arr = New-Object 'object[,]' 1,21
$arr[0,0] = 'Tutaenui Stream at d/s Marton STP'
$arr[0,1] = 'Total Phosphorus (95th Percentile)'
$arr[0,3] = '2019 - 2023'
$arr[0,4] = 'Impact'
$arr[0,5] = 1.24
$arr[0,6] = 2.19681034482759
$arr[0,7] = 7.15
$arr[0,8] = 6.102
$arr[0,11] = 3.53
$arr[0,12] = 4.9112
$arr[0,13] = 5.4244
$arr[0,14] = 1803578.705
$arr[0,15] = 5557699.998
$arr[0,16] = 'Rangitikei District'
$arr[0,17] = 'Rangitīkei-Turakina'
$arr[0,18] = 'Coastal Rangitikei'
$arr[0,19] = 'Rang_4d'
$arr[0,20] = 'g/m3'
$ws.Range("A181:U181").Value = $arr

$arr = New-Object 'object[,]' 1,21
$arr[0,0] = 'Tutaenui Stream at d/s Marton STP'
$arr[0,1] = 'Total Phosphorus (Median)'
$arr[0,3] = '2019 - 2023'
$arr[0,4] = 'Impact'
$arr[0,5] = 1.24
$arr[0,6] = 2.19681034482759
$arr[0,7] = 7.15
$arr[0,8] = 6.102
$arr[0,11] = 3.53
$arr[0,12] = 4.9112
$arr[0,13] = 5.4244
$arr[0,14] = 1803578.705
$arr[0,15] = 5557699.998
$arr[0,16] = 'Rangitikei District'
$arr[0,17] = 'Rangitīkei-Turakina'
$arr[0,18] = 'Coastal Rangitikei'
$arr[0,19] = 'Rang_4d'
$arr[0,20] = 'g/m3'
$ws.Range("A182:U182").Value = $arr

